$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Start"
$ws.Range("B2").Value = "IcecreamShop,0.5,5,0"
$ws.Range("C2").Value = "Subway St.1,3,10,5"

# Row 3
$ws.Range("A3").Value = "Subway St.1"
$ws.Range("B3").Value = "Subway St.2,25,35,15"
$ws.Range("C3").Value = "Subway St.3,15,20,10"
$ws.Range("D3").Value = "Start,3,10,5"

# Row 4
$ws.Range("A4").Value = "Subway St.2"
$ws.Range("B4").Value = "mall,0.3,5,0"
$ws.Range("C4").Value = "Subway St.1,25,35,15"

# Row 5
$ws.Range("A5").Value = "Subway St.3"
$ws.Range("B5").Value = "Subway St.1,15,20,10"
$ws.Range("C5").Value = "Cafe,5,10,10"
$ws.Range("D5").Value = "Aquarium,3,7,5"
$ws.Range("E5").Value = "GameCenter,8,10,20"

# Row 6
$ws.Range("A6").Value = "Cafe"
$ws.Range("B6").Value = "Subway St.3,5,10,10"
$ws.Range("C6").Value = "Aquarium,1000,120,20000"

# Row 7
$ws.Range("A7").Value = "Aquarium"
$ws.Range("B7").Value = "Cafe,1000,120,20000"
$ws.Range("C7").Value = "Subway St.3,3,7,5"
$ws.Range("D7").Value = "GameCenter,4,5,2"

# Row 8
$ws.Range("A8").Value = "GameCenter"
$ws.Range("B8").Value = "Aquarium,4,5,2"
$ws.Range("C8").Value = "Subway St.3,8,10,20"

# Row 9
$ws.Range("A9").Value = "IcecreamShop"
$ws.Range("B9").Value = "Start,0.5,5,0"

# Row 10
$ws.Range("A10").Value = "mall"
$ws.Range("B10").Value = "Subway St.2,0.3,5,0"

# Row 11 - blank but styled cells (no values)

# Apply the same center-aligned cell style ("s=1") used throughout the sheet
# to every cell in the A2:E8 block (rows 2-8 use cols A:E), rows 9-10
# (cols A:C only) and row 11 (cols A:B only), including the new, empty
# cells, so that all cells carry consistent formatting as in the source
# workbook.
$ws.Range("A2:E8").HorizontalAlignment = -4108
$ws.Range("A9:C10").HorizontalAlignment = -4108
$ws.Range("A11:B11").HorizontalAlignment = -4108

# Column widths
$ws.Columns("A:K").ColumnWidth = 33.6

# Selection
$ws.Range("D15").Select()
